# Generate Report for Handoff
# Replaces the old GUID-named handoff artifact (c4604868-4387-4522-a64b-39910525caec)
# with the newly generated one (80b866fe-2773-43fc-a113-ee22b175f378) across the
# Overview / zh-cn / de-de sheets, and refreshes the associated handoff/handback
# timestamps.

$wb = $excel.ActiveWorkbook

$oldGuid = "c4604868-4387-4522-a64b-39910525caec"
$newGuid = "80b866fe-2773-43fc-a113-ee22b175f378"

$hyperlinkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ce8d00323b827a84e8dc1915cef0197ad5ce2277/e2e/$oldGuid.md"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newGuid.md"

$overviewDisplay = "e2e\$newGuid.md"
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $hyperlinkAddress, "", "", $overviewDisplay) | Out-Null

$wsOverview.Range("G2").Value = "2016-08-17 00:53:52"

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$zhDisplay = "$newGuid.md"
$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $hyperlinkAddress, "", "", $zhDisplay) | Out-Null

$wsZhCn.Range("G2").Value = "$newGuid.f12103602e2d0942812f0a07d8cc556a1cec5fcd.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-17 00:53:47"

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$deDisplay = "$newGuid.md"
$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $hyperlinkAddress, "", "", $deDisplay) | Out-Null

$wsDeDe.Range("G2").Value = "$newGuid.f12103602e2d0942812f0a07d8cc556a1cec5fcd.de-de.xlf"

Write-Host "Handoff report regenerated: $oldGuid -> $newGuid"
